$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Structural: insert a new blank row above row 1 (shifts header/data
# down by one row) and empty out the old index column (A).
# ------------------------------------------------------------------
$ws.Rows.Item(1).Insert()
$ws.Columns.Item(1).Clear()

# ------------------------------------------------------------------
# Column widths / row heights
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15.83203125
$ws.Columns.Item(2).ColumnWidth = 15.83203125
$ws.Range($ws.Columns.Item(3), $ws.Columns.Item(16384)).ColumnWidth = 15.83203125

$ws.Range("1:8").RowHeight = 27

# ------------------------------------------------------------------
# Helper groups of cells, matching the final layout:
#   row1          -> blank spacer row                (bold style, bottom-medium border)
#   row2 (B2:J2)  -> header                           (bold style, bottom-medium border)
#   row8 (B8:J8)  -> footer                           (B8 bold+border / C8:J8 regular+border)
#   row3 B3, row5 B5          -> bold, no border
#   row3 C3:J3, row5 C5:J5    -> regular, no border
#   row4 B4, row6 B6, row7 B7 -> bold, top+bottom thin border
#   row4 C4:J4, row6 C6:J6, row7 C7:J7 -> regular, top+bottom thin border
# ------------------------------------------------------------------

function Format-Range($r, [bool]$bold, [int]$borderKind) {
    # borderKind: 0 = none, 1 = bottom medium only, 2 = top+bottom thin
    $r.Font.Name = "Cambria Math"
    $r.Font.Size = 11
    $r.Font.Bold = $bold
    $r.Interior.Pattern = 1
    $r.Interior.ThemeColor = 2
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4108
    if ($borderKind -eq 1) {
        $r.Borders.Item(9).LineStyle = 1
        $r.Borders.Item(9).Weight = -4138
    } elseif ($borderKind -eq 2) {
        $r.Borders.Item(8).LineStyle = 1
        $r.Borders.Item(8).Weight = 2
        $r.Borders.Item(9).LineStyle = 1
        $r.Borders.Item(9).Weight = 2
    }
}

# Blank spacer row 1 (B1 only has a cell)
Format-Range $ws.Range("B1") $true 1

# Header row 2
Format-Range $ws.Range("B2:J2") $true 1

# Data row 3 (regular)
Format-Range $ws.Range("B3") $true 0
Format-Range $ws.Range("C3:J3") $false 0

# Data row 4 (customFormat / bold col B)
Format-Range $ws.Range("B4") $true 2
Format-Range $ws.Range("C4:J4") $false 2

# Data row 5 (regular)
Format-Range $ws.Range("B5") $true 0
Format-Range $ws.Range("C5:J5") $false 0

# Data row 6 (customFormat / bold col B)
Format-Range $ws.Range("B6") $true 2
Format-Range $ws.Range("C6:J6") $false 2

# Data row 7 (customFormat / bold col B)
Format-Range $ws.Range("B7") $true 2
Format-Range $ws.Range("C7:J7") $false 2

# Footer row 8
Format-Range $ws.Range("B8") $true 1
Format-Range $ws.Range("C8:J8") $false 1

# ------------------------------------------------------------------
# Selection matches the saved view in the target file
# ------------------------------------------------------------------
$ws.Range("C10").Select()
